$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mBom")

$ws.Range("F37").Value = "CONN HEADER .050"" 14PS DL PCB AU"
$ws.Range("H37").Value = "GRPB072VWVN-RC"
$ws.Range("I37").Value = "S9015E-07-ND"

$ws.Range("K37").Value = 1.26
$ws.Range("L37").Value = 1.048
$ws.Range("M37").Value = 0.776
$ws.Range("N37").Value = 0.5044
